$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")
$ws.Activate()

# "there will be at this moment 4 new mod setups" - rescale the mod-setup
# bounds used to derive the H column (ALT L HP -> ALT H HP range).
# C2: 3 -> 0 , D2: 36 -> 30
# Every H-column cell (H20:H164) holds
#   =ROUND(((F-$C$4)/($D$4-$C$4))*($D$2-$C$2),0)+$C$2
# so they all recompute automatically off these two inputs.
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 30

# Move the selection/cursor to D3 (was C17); the stale scroll position
# (topLeftCell="A130") is cleared as part of the same view refresh.
$ws.Range("D3").Select() | Out-Null
